# Append new player stat rows (14-28) to the TE aggregate sheet,
# continuing the existing alternating row-style banding.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style source ranges already present on the sheet (A2:E4 = style "2",
# A5:E7 = style "3"). Copy them onto the new row blocks first so the
# banding continues, then overwrite the values/labels.
$styleA = $ws.Range("A2:E4")
$styleB = $ws.Range("A5:E7")

$styleA.Copy($ws.Range("A14:E16"))
$styleB.Copy($ws.Range("A17:E19"))
$styleA.Copy($ws.Range("A20:E22"))
$styleB.Copy($ws.Range("A23:E25"))
$styleA.Copy($ws.Range("A26:E28"))

# New data rows (player, season group, Y/R, Y/Tgt, Succ%)
$data = @(
    @(14, "John Mundt",   "Group1",     9.6,               8.666666666666666,   80.96666666666667),
    @(15, "John Mundt",   "Group2",     8.333333333333334, 6.5,                 54.43333333333334),
    @(16, "John Mundt",   "Difference", -1.266666666666666,-2.166666666666666, -26.53333333333333),

    @(17, "Mike Gesicki", "Group1",     11.73333333333333, 7.233333333333334,   52.03333333333333),
    @(18, "Mike Gesicki", "Group2",     9.966666666666667, 6.8,                 51.73333333333333),
    @(19, "Mike Gesicki", "Difference", -1.766666666666667,-0.4333333333333345,-0.3000000000000043),

    @(20, "Travis Kelce", "Group1",     12.8,               9.066666666666668,   63.53333333333333),
    @(21, "Travis Kelce", "Group2",     10.43333333333333, 7.699999999999999,   61.8),
    @(22, "Travis Kelce", "Difference", -2.366666666666667,-1.366666666666669, -1.733333333333327),

    @(23, "Nick Vannett",  "Group1",     10.43333333333333, 6.977777777777779,   52.9),
    @(24, "Nick Vannett",  "Group2",     6.544444444444444, 5.333333333333333,   80.3888888888889),
    @(25, "Nick Vannett",  "Difference", -3.888888888888889,-1.644444444444446,  27.4888888888889),

    @(26, "Noah Fant",     "Group1",     11.63333333333333, 7.7,                 47.26666666666667),
    @(27, "Noah Fant",     "Group2",     11,                8.366666666666667,   57.43333333333334),
    @(28, "Noah Fant",     "Difference", -0.6333333333333329,0.666666666666667,  10.16666666666666)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
